$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old columns beyond the new data extent (columns E:I)
$ws.Range("E1:I2").EntireColumn.Delete()

# Header row
$ws.Range("A1").Value = "id_cliente"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "apellido"
$ws.Range("D1").Value = "cedula"

# Data row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Andres"
$ws.Range("C2").Value = "guido"
$ws.Range("D2").Value = "'123"
$ws.Range("D2").Style = "Normal"
